$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous last row (42) loses its "last row" date formatting and
# becomes a regular data row; the new row 43 becomes the new "last row".
$ws.Range("A42").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A43").Value = 45628
$ws.Range("A43").NumberFormat = "YYYY-MM-DD"
$ws.Range("B43").Value = 112
$ws.Range("C43").Value = 93
$ws.Range("D43").Value = 102
